$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Npc")
$ws2 = $wb.Worksheets.Item("디스크립션")

# ============================================================
# Sheet "Npc": header rename + skill data for 5 monsters
# (No.107, No.108, Jake_B, Jake_R, Doncina)
# ============================================================

# --- Header row ---
$ws.Range("C1").Value = "nameKor"
$ws.Range("H1").Value = "recognizeValue"
$ws.Range("I1").Value = "skillValue"
$ws.Range("J1").Value = "dropTable"

# --- Row 2 (No.107) ---
$ws.Range("A2").Value = 19100
$ws.Range("B2").Value = "no107_01"
$ws.Range("C2").Value = "넘버107"
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = "Mob"
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = "{(19_DeadlyAttack)}"
$ws.Range("J2").Value = "{(1100,5),(2100,5),(3100,5)}"

# --- Row 3 (No.108) ---
$ws.Range("A3").Value = 19101
$ws.Range("B3").Value = "no107_02"
$ws.Range("C3").Value = "넘버108"
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = "Mob"
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = "{(19_CriticalHit)}"
$ws.Range("J3").Value = "{(1100,5),(2101,5),(3100,5)}"

# --- Row 4 (Jake_B, was the old "Jake" row) ---
$ws.Range("A4").Value = 19102
$ws.Range("B4").Value = "Jake_B"
$ws.Range("C4").Value = "썩은쥐돌이"
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = "Mob"
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = "{(19_jakebounce)}"
$ws.Range("J4").Value = "{(1102,5),(2102,5),(3101,5)}"

# --- Row 5 (jake_R, new row) ---
$ws.Range("A5").Value = 19103
$ws.Range("B5").Value = "jake_R"
$ws.Range("C5").Value = "엉성한쥐돌이"
$ws.Range("E5").Value = 30
$ws.Range("F5").Value = "Mob"
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = "{(19_jakerange)}"
$ws.Range("J5").Value = "{(1101,5),(2101,5),(3102,5)}"

# --- Row 6 (doncina, new row) ---
$ws.Range("A6").Value = 19200
$ws.Range("B6").Value = "doncina"
$ws.Range("C6").Value = "돈시나"
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = "Boss"
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = "{(19_MustleMustle),(19_HustleHustle)}"
$ws.Range("J6").Value = "{(1103,5),(2100,5),(3101,5),(1101,5),(2101,5),(3100,5),(4101,5)}"

# Drop-table column (J) is no longer highlighted in yellow -- clear the fill
$ws.Range("J1:J6").ClearFormats() | Out-Null
$ws.Range("J1").Value = "dropTable"
$ws.Range("J2").Value = "{(1100,5),(2100,5),(3100,5)}"
$ws.Range("J3").Value = "{(1100,5),(2101,5),(3100,5)}"
$ws.Range("J4").Value = "{(1102,5),(2102,5),(3101,5)}"
$ws.Range("J5").Value = "{(1101,5),(2101,5),(3102,5)}"
$ws.Range("J6").Value = "{(1103,5),(2100,5),(3101,5),(1101,5),(2101,5),(3100,5),(4101,5)}"

# --- Column widths (bestFit) ---
$ws.Columns.Item(3).ColumnWidth = 12.2857142857143
$ws.Columns.Item(8).ColumnWidth = 14.4285714285714
$ws.Columns.Item(9).ColumnWidth = 35.8571428571429
$ws.Columns.Item(10).ColumnWidth = 47.4285714285714

# --- View / selection state ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws.Range("D14").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("D4").Select() | Out-Null

$ws.Activate() | Out-Null
